# Refresh the 'cryptos' price/volume snapshot (Price column D, Volume(1h) column E)
# on Sheet1. Every cell written here already holds a text value (t="inlineStr")
# in the source workbook, so we re-write each one as text too:
#  - Values that don't look like a plain number (e.g. thousands-grouped prices
#    like '62.748.82', percents like '  -0.55%  ', or subscript-digit prices
#    like '0.0₃0954') round-trip as text on their own.
#  - Values that DO look like a plain number (e.g. '570.90', '1.70') would
#    otherwise be auto-coerced to a numeric cell by Excel (dropping trailing
#    zeros, etc.), so a leading apostrophe forces text entry; the cell style is
#    then reset to Normal so we don't leave a stray 'quote prefix' format behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.748.82'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = '2.457.87'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''570.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.24%  '
$ws.Range('D6').Value = '''146.23'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -1.91%  '
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('D10').Value = '''0.162'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.21%  '
$ws.Range('E11').Value = '  -2.26%  '
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('D13').Value = '''28.49'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.91%  '
$ws.Range('E14').Value = '  -3.30%  '
$ws.Range('D15').Value = '2.900.70'
$ws.Range('E15').Value = '  -0.74%  '
$ws.Range('D16').Value = '62.573.79'
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('D17').Value = '2.452.81'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').Value = '''7.66'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.70%  '
$ws.Range('D19').Value = '''10.72'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E20').Value = '  -0.77%  '
$ws.Range('D21').Value = '''321.45'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.32%  '
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E24').Value = '  +2.30%  '
$ws.Range('D25').Value = '''64.75'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.39%  '
$ws.Range('D26').Value = '''645.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.37%  '
$ws.Range('D28').Value = '0.0₃0954'
$ws.Range('E28').Value = '  -3.87%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('E30').Value = '  -3.48%  '
$ws.Range('D31').Value = '''7.83'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.73%  '
$ws.Range('D32').Value = '''1.81'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.98%  '
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('E35').Value = '  -4.15%  '
$ws.Range('E36').Value = '  -3.39%  '
$ws.Range('D37').Value = '''150.34'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.05%  '
$ws.Range('D38').Value = '''18.52'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.52%  '
$ws.Range('E39').Value = '  -2.40%  '
$ws.Range('E40').Value = '  -3.01%  '
$ws.Range('D41').Value = '''2.64'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.36%  '
$ws.Range('D42').Value = '''1.70'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.65%  '
$ws.Range('D43').Value = '0.0₆0315'
$ws.Range('E43').Value = '  +2.91%  '
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('D45').Value = '''152.32'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('D46').Value = '''15.41'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.87%  '
$ws.Range('D47').Value = '''3.52'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.33%  '
$ws.Range('D48').Value = '''0.603'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('D49').Value = '''19.98'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.35%  '
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('D51').Value = '''0.0901'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.95%  '
